$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column F (6) so the new "avance" notes fit; matches stored width ~40.83
$ws.Columns.Item(6).ColumnWidth = 40

# Highlight rows 2, 3, 6 and 8 (the activities marked as done) with the
# accent6 theme color fill, keeping the existing center+wrap alignment.
function Set-DoneRow($range) {
    $range.HorizontalAlignment = -4108
    $range.WrapText = $true
    $range.Interior.ThemeColor = 10
}

Set-DoneRow($ws.Range("A2:F2"))
Set-DoneRow($ws.Range("A3:F3"))
Set-DoneRow($ws.Range("A6:F6"))
Set-DoneRow($ws.Range("A8:F8"))

# Move the active selection to the row with the new activity (row 9)
$ws.Range("A9:F9").Select()
